# Generate Report for Handoff
# Adds a new localization-status row (for file
# "45e61196-dd2b-4863-a39e-77c67d4820a6") to the Overview, zh-cn and
# de-de worksheets, mirroring the existing row 2 for
# "304848e7-1a06-40e6-bfe6-b511954e83f4".

$wb = $excel.ActiveWorkbook

# Hyperlink font look-alike constants (underline + the workbook's custom
# "HyperLink" blue FF6495ED, expressed as a BGR COM color value).
$hlUnderline = -4142
$hlColor = 15570276

function Style-AsHyperlink($rng) {
    $rng.Font.Underline = $hlUnderline
    $rng.Font.Color = $hlColor
    $rng.Font.Name = "Calibri"
    $rng.Font.Size = 11
}

# ---------------------------------------------------------------------
# Overview sheet
# ---------------------------------------------------------------------
$wsOverview = $wb.Worksheets.Item("Overview")

$wsOverview.Rows("2:2").Copy()
$wsOverview.Rows("3:3").Insert(-4121)

$wsOverview.Range("A3").Value = "45e61196-dd2b-4863-a39e-77c67d4820a6.md"
$wsOverview.Range("B3").Value = "Ready for handoff"
$wsOverview.Range("C3").Value = "Ready for handoff"
$wsOverview.Range("D3").Value = "2016-25-11 22:25:38"

$wsOverview.Hyperlinks.Add(
    $wsOverview.Range("A3"),
    "https://github.com/OpenLocalizationTest/oltest/blob/1d6a735efd6d549edb8ff954fd2ba5dbe631ee69/e2e/45e61196-dd2b-4863-a39e-77c67d4820a6.md",
    "",
    "",
    "45e61196-dd2b-4863-a39e-77c67d4820a6.md")
Style-AsHyperlink $wsOverview.Range("A3")

# ---------------------------------------------------------------------
# zh-cn sheet
# ---------------------------------------------------------------------
$wsZhCn = $wb.Worksheets.Item("zh-cn")

$wsZhCn.Rows("2:2").Copy()
$wsZhCn.Rows("3:3").Insert(-4121)

$wsZhCn.Range("A3").Value = "45e61196-dd2b-4863-a39e-77c67d4820a6.md"
$wsZhCn.Range("B3").Value = ".md"
$wsZhCn.Range("C3").Value = "Ready for handoff"
$wsZhCn.Range("D3").Value = "45e61196-dd2b-4863-a39e-77c67d4820a6.f55201f77774808f6ba27d79737dc1a0cecfa6a3.zh-cn.xlf"
$wsZhCn.Range("E3").Value = "2016-03-11 22:25:36"
$wsZhCn.Range("H3").Value = "0001-01-01 00:00:00"
$wsZhCn.Range("I3").Value = "Include"

$wsZhCn.Hyperlinks.Add(
    $wsZhCn.Range("A3"),
    "https://github.com/OpenLocalizationTest/oltest/blob/1d6a735efd6d549edb8ff954fd2ba5dbe631ee69/e2e/45e61196-dd2b-4863-a39e-77c67d4820a6.md",
    "",
    "",
    "45e61196-dd2b-4863-a39e-77c67d4820a6.md")
Style-AsHyperlink $wsZhCn.Range("A3")

$wsZhCn.Hyperlinks.Add(
    $wsZhCn.Range("B3"),
    "https://github.com/OpenLocalizationTest/oltest/blob/1d6a735efd6d549edb8ff954fd2ba5dbe631ee69/e2e/45e61196-dd2b-4863-a39e-77c67d4820a6.md",
    "",
    "",
    ".md")
Style-AsHyperlink $wsZhCn.Range("B3")

$wsZhCn.Hyperlinks.Add(
    $wsZhCn.Range("D3"),
    "https://github.com/OpenLocalizationTestOrg/olhandoff/blob/c61ea852603dcfc67c917c8b051a2d941426a724/ol-handoff/OpenLocalizationTestOrg/oltest.zh-cn/ci/ht/45e61196-dd2b-4863-a39e-77c67d4820a6.f55201f77774808f6ba27d79737dc1a0cecfa6a3.zh-cn.xlf",
    "",
    "",
    "45e61196-dd2b-4863-a39e-77c67d4820a6.f55201f77774808f6ba27d79737dc1a0cecfa6a3.zh-cn.xlf")
Style-AsHyperlink $wsZhCn.Range("D3")

# ---------------------------------------------------------------------
# de-de sheet
# ---------------------------------------------------------------------
$wsDeDe = $wb.Worksheets.Item("de-de")

$wsDeDe.Rows("2:2").Copy()
$wsDeDe.Rows("3:3").Insert(-4121)

$wsDeDe.Range("A3").Value = "45e61196-dd2b-4863-a39e-77c67d4820a6.md"
$wsDeDe.Range("B3").Value = ".md"
$wsDeDe.Range("C3").Value = "Ready for handoff"
$wsDeDe.Range("D3").Value = "45e61196-dd2b-4863-a39e-77c67d4820a6.f55201f77774808f6ba27d79737dc1a0cecfa6a3.de-de.xlf"
$wsDeDe.Range("E3").Value = "2016-03-11 22:25:38"
$wsDeDe.Range("H3").Value = "0001-01-01 00:00:00"
$wsDeDe.Range("I3").Value = "Include"

$wsDeDe.Hyperlinks.Add(
    $wsDeDe.Range("A3"),
    "https://github.com/OpenLocalizationTest/oltest/blob/1d6a735efd6d549edb8ff954fd2ba5dbe631ee69/e2e/45e61196-dd2b-4863-a39e-77c67d4820a6.md",
    "",
    "",
    "45e61196-dd2b-4863-a39e-77c67d4820a6.md")
Style-AsHyperlink $wsDeDe.Range("A3")

$wsDeDe.Hyperlinks.Add(
    $wsDeDe.Range("B3"),
    "https://github.com/OpenLocalizationTest/oltest/blob/1d6a735efd6d549edb8ff954fd2ba5dbe631ee69/e2e/45e61196-dd2b-4863-a39e-77c67d4820a6.md",
    "",
    "",
    ".md")
Style-AsHyperlink $wsDeDe.Range("B3")

$wsDeDe.Hyperlinks.Add(
    $wsDeDe.Range("D3"),
    "https://github.com/OpenLocalizationTestOrg/olhandoff/blob/f34ab406754316f5835a23911e15d26aab3807df/ol-handoff/OpenLocalizationTestOrg/oltest.de-de/ci/ht/45e61196-dd2b-4863-a39e-77c67d4820a6.f55201f77774808f6ba27d79737dc1a0cecfa6a3.de-de.xlf",
    "",
    "",
    "45e61196-dd2b-4863-a39e-77c67d4820a6.f55201f77774808f6ba27d79737dc1a0cecfa6a3.de-de.xlf")
Style-AsHyperlink $wsDeDe.Range("D3")

Write-Output "Added handoff row for 45e61196-dd2b-4863-a39e-77c67d4820a6 to Overview, zh-cn, de-de"
